# Automated tracker sync: fill in results/profit for events that have now
# been settled, and append the newest tracked event as a new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Resolved bets: set resultado (G) / profit (H) ---------------------
$updates = @(
    @{ Row = 32; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 58; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 61; Resultado = "Acierto"; Profit = 2.4 },
    @{ Row = 65; Resultado = "Acierto"; Profit = 1.5 },
    @{ Row = 69; Resultado = "Acierto"; Profit = 1.75 },
    @{ Row = 70; Resultado = "Acierto"; Profit = 1.2 },
    @{ Row = 74; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 75; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 76; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 78; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 84; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 85; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 86; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 87; Resultado = "Fallo";   Profit = -1 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Resultado
    $ws.Cells.Item($u.Row, 8).Value = $u.Profit
}

# --- 2. Append the newest event as row 88 ----------------------------------
$newRow = 88

$ws.Cells.Item($newRow, 1).Value = 14580798

# "fecha" looks like a date ("2025-09-02"); entering it as a plain Value
# makes Excel auto-convert it to a date serial, so build it as a text
# formula first and then harden it down to a literal value/string.
$ws.Cells.Item($newRow, 2).Formula = '="2025-09-02"'
$ws.Cells.Item($newRow, 2).Copy()
$ws.Cells.Item($newRow, 2).PasteSpecial(-4163)

$ws.Cells.Item($newRow, 3).Value = "Thiago Monteiro"
$ws.Cells.Item($newRow, 4).Value = "Francesco Passaro"
$ws.Cells.Item($newRow, 5).Value = "Gana Thiago Monteiro"
$ws.Cells.Item($newRow, 6).Value = 3

# Outcome is still pending for this event, so resultado/profit stay blank
# (matching every other not-yet-settled row in the sheet).
$ws.Cells.Item($newRow, 7).Formula = '=""'
$ws.Cells.Item($newRow, 8).Formula = '=""'
